# Atualização de bases das ligas, do dia: 01-05-2024 às 16:15
#
# The underlying data table (rows 2-139, columns B:AB) is re-sorted by the
# "Date" column (D) in ascending order. The leading index column (A) is a
# fixed 0-based row counter and must stay untouched. Because the dataset is
# already sorted almost everywhere, this resolves to swapping the full
# content (columns B through AB) of seven pairs of adjacent rows whose
# timestamps were out of order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB")

$pairs = @(
    @(3,4),
    @(16,17),
    @(19,20),
    @(36,37),
    @(55,56),
    @(68,69),
    @(114,115)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $vals1 = @{}
    $vals2 = @{}

    foreach ($c in $cols) {
        $vals1[$c] = $ws.Range("$c$r1").Value2
        $vals2[$c] = $ws.Range("$c$r2").Value2
    }

    foreach ($c in $cols) {
        $ws.Range("$c$r1").Value = $vals2[$c]
        $ws.Range("$c$r2").Value = $vals1[$c]
    }
}

Write-Host "Swapped $($pairs.Count) row pairs to restore ascending Date order."
